# Update NATMI LR-pair results (Bmp3-Acvr2b) with recomputed values based on
# the new per-cell TPM input. Only the receptor-side statistics (columns
# G/H ligand avg/total - negligible float re-summation noise - and
# K/L/M/N receptor-expressing-cells/detection-rate/avg/total expression,
# plus the derived specificity columns O/P/Q/R/S/T that are recalculated
# from them) change; the ligand/receptor identity columns (A-F, I, J) are
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3037726666666666
$ws.Range("H2").Value = 0.911318
$ws.Range("M2").Value = 0.978762
$ws.Range("N2").Value = 2.936286
$ws.Range("O2").Value = 0.3819465121442868
$ws.Range("P2").Value = 0.3819465121442868
$ws.Range("Q2").Value = 0.297321142772
$ws.Range("R2").Value = 2.675890284948
$ws.Range("S2").Value = 0.3708789665719145
$ws.Range("T2").Value = 0.3708789665719145

# Row 3
$ws.Range("G3").Value = 0.3037726666666666
$ws.Range("H3").Value = 0.911318
$ws.Range("O3").Value = 0.0406458950147437
$ws.Range("P3").Value = 0.04064589501474371
$ws.Range("S3").Value = 0.03946811152647454
$ws.Range("T3").Value = 0.03946811152647455

# Row 4
$ws.Range("G4").Value = 0.3037726666666666
$ws.Range("H4").Value = 0.911318
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3139526666666667
$ws.Range("N4").Value = 0.9418580000000001
$ws.Range("O4").Value = 0.1225151017425393
$ws.Range("P4").Value = 0.1225151017425393
$ws.Range("Q4").Value = 0.09537023876044445
$ws.Range("R4").Value = 0.858332148844
$ws.Range("S4").Value = 0.1189650196532253
$ws.Range("T4").Value = 0.1189650196532253

# Row 5
$ws.Range("G5").Value = 0.3037726666666666
$ws.Range("H5").Value = 0.911318
$ws.Range("M5").Value = 0.5829876666666666
$ws.Range("N5").Value = 1.748963
$ws.Range("O5").Value = 0.2275017888991087
$ws.Range("P5").Value = 0.2275017888991087
$ws.Range("Q5").Value = 0.1770957181371111
$ws.Range("R5").Value = 1.593861463234
$ws.Range("S5").Value = 0.220909540151237
$ws.Range("T5").Value = 0.2209095401512371

# Row 6
$ws.Range("G6").Value = 0.3037726666666666
$ws.Range("H6").Value = 0.911318
$ws.Range("M6").Value = 0.3019996666666667
$ws.Range("N6").Value = 0.905999
$ws.Range("O6").Value = 0.117850631054404
$ws.Range("P6").Value = 0.117850631054404
$ws.Range("Q6").Value = 0.09173924407577777
$ws.Range("R6").Value = 0.825653196682
$ws.Range("S6").Value = 0.1144357098849322
$ws.Range("T6").Value = 0.1144357098849322

# Row 7
$ws.Range("G7").Value = 0.3037726666666666
$ws.Range("H7").Value = 0.911318
$ws.Range("M7").Value = 0.2807033333333333
$ws.Range("N7").Value = 0.8421099999999999
$ws.Range("O7").Value = 0.1095400711449175
$ws.Range("P7").Value = 0.1095400711449176
$ws.Range("Q7").Value = 0.08527000010888887
$ws.Range("R7").Value = 0.7674300009799999
$ws.Range("S7").Value = 0.1063659624913496
$ws.Range("T7").Value = 0.1063659624913496

# Row 8
$ws.Range("M8").Value = 0.978762
$ws.Range("N8").Value = 2.936286
$ws.Range("O8").Value = 0.3819465121442868
$ws.Range("P8").Value = 0.3819465121442868
$ws.Range("Q8").Value = 0.008872477530000001
$ws.Range("R8").Value = 0.07985229777
$ws.Range("S8").Value = 0.01106754557237234
$ws.Range("T8").Value = 0.01106754557237234

# Row 9
$ws.Range("O9").Value = 0.0406458950147437
$ws.Range("P9").Value = 0.04064589501474371
$ws.Range("S9").Value = 0.001177783488269161
$ws.Range("T9").Value = 0.001177783488269161

# Row 10
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3139526666666667
$ws.Range("N10").Value = 0.9418580000000001
$ws.Range("O10").Value = 0.1225151017425393
$ws.Range("P10").Value = 0.1225151017425393
$ws.Range("Q10").Value = 0.002845980923333334
$ws.Range("R10").Value = 0.02561382831
$ws.Range("S10").Value = 0.003550082089314007
$ws.Range("T10").Value = 0.003550082089314007

# Row 11
$ws.Range("M11").Value = 0.5829876666666666
$ws.Range("N11").Value = 1.748963
$ws.Range("O11").Value = 0.2275017888991087
$ws.Range("P11").Value = 0.2275017888991087
$ws.Range("Q11").Value = 0.005284783198333333
$ws.Range("R11").Value = 0.047563048785
$ws.Range("S11").Value = 0.006592248747871644
$ws.Range("T11").Value = 0.006592248747871646

# Row 12
$ws.Range("M12").Value = 0.3019996666666667
$ws.Range("N12").Value = 0.905999
$ws.Range("O12").Value = 0.117850631054404
$ws.Range("P12").Value = 0.117850631054404
$ws.Range("Q12").Value = 0.002737626978333334
$ws.Range("R12").Value = 0.024638642805
$ws.Range("S12").Value = 0.003414921169471832
$ws.Range("T12").Value = 0.003414921169471832

# Row 13
$ws.Range("M13").Value = 0.2807033333333333
$ws.Range("N13").Value = 0.8421099999999999
$ws.Range("O13").Value = 0.1095400711449175
$ws.Range("P13").Value = 0.1095400711449176
$ws.Range("Q13").Value = 0.002544575716666666
$ws.Range("R13").Value = 0.02290118145
$ws.Range("S13").Value = 0.003174108653567966
$ws.Range("T13").Value = 0.003174108653567967
